$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9858714938163757
$ws.Range("B1").Value = 1.773585915565491
$ws.Range("C1").Value = 4.945278167724609
$ws.Range("D1").Value = 1.323824167251587
$ws.Range("E1").Value = 1.29129421710968
